$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that must stay as Text even when it parses as a
# plain number (e.g. "268.67"); Column E (percent strings) and B/C (names/links)
# never look like pure numbers so a plain .Value assignment keeps them as text.
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

$ws.Range('D2').Value = '57.222.12'
$ws.Range('E2').Value = '  +1.30%  '
$ws.Range('D3').Value = '2.428.08'
$ws.Range('E3').Value = '  -2.22%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '489.79'
$ws.Range('E5').Value = '  -0.28%  '
Set-TextValue 'D6' '155.87'
$ws.Range('E6').Value = '  +4.04%  '
Set-TextValue 'D7' '0.621'
$ws.Range('E7').Value = '  +20.41%  '
Set-TextValue 'D8' '0.996'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').Value = '2.450.22'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  +2.15%  '
Set-TextValue 'D11' '5.68'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('D14').Value = '2.852.67'
$ws.Range('E14').Value = '  -2.12%  '
$ws.Range('D15').Value = '57.235.46'
$ws.Range('E15').Value = '  +1.35%  '
Set-TextValue 'D16' '20.87'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('E17').Value = '  -1.35%  '
$ws.Range('D18').Value = '2.443.01'
$ws.Range('E18').Value = '  -1.91%  '
Set-TextValue 'D19' '4.78'
$ws.Range('E19').Value = '  +5.71%  '
Set-TextValue 'D20' '329.64'
$ws.Range('E20').Value = '  +2.78%  '
Set-TextValue 'D21' '10.00'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  +0.49%  '
Set-TextValue 'D24' '58.31'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('E25').Value = '  +0.18%  '
Set-TextValue 'D26' '0.998'
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').Value = '2.539.54'
$ws.Range('E28').Value = '  -2.23%  '
Set-TextValue 'D29' '7.33'
$ws.Range('E29').Value = '  -4.08%  '
$ws.Range('E30').Value = '  -0.17%  '
Set-TextValue 'D31' '0.998'
$ws.Range('E31').Value = '  -0.02%  '
Set-TextValue 'D32' '18.80'
$ws.Range('E32').Value = '  +2.49%  '
Set-TextValue 'D33' '149.34'
$ws.Range('E33').Value = '  -0.60%  '
Set-TextValue 'D34' '1.53'
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('E35').Value = '  +2.59%  '
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('E37').Value = '  -1.07%  '
Set-TextValue 'D38' '0.860'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('E39').Value = '  +10.51%  '
Set-TextValue 'D40' '34.26'
$ws.Range('E40').Value = '  +1.15%  '
$ws.Range('E41').Value = '  -0.65%  '
Set-TextValue 'D42' '3.54'
$ws.Range('E42').Value = '  +1.22%  '
Set-TextValue 'D43' '0.995'
$ws.Range('E43').Value = '  -0.08%  '
Set-TextValue 'D44' '0.599'
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D45' '268.67'
$ws.Range('E45').Value = '  +1.43%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D46' '0.0537'
$ws.Range('E46').Value = '  -3.74%  '
Set-TextValue 'D47' '0.0229'
$ws.Range('E47').Value = '  -0.36%  '
Set-TextValue 'D48' '10.22'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('E49').Value = '  -3.18%  '
Set-TextValue 'D50' '17.59'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').Value = '1.861.56'
$ws.Range('E51').Value = '  -1.68%  '
